# Added execute command utility method
# Split the combined "logindata" sheet into three sheets:
#   - logindata   : keeps only the login test rows (A1:C5)
#   - PMTestData  : the "create/edit user" rows that used to live at the
#                   bottom of logindata (rows 6-7), now its own sheet
#   - SNMTestData : a full copy of the original logindata content (A1:C7)
# and leaves the existing "IP" sheet as the last tab.

$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("logindata")

# --- PMTestData -----------------------------------------------------
# Full copy of logindata placed right after it, then trim down to just
# the header + the two "create/edit user" rows (old rows 6 & 7).
$login.Copy($null, $login) | Out-Null
$pm = $wb.ActiveSheet
$pm.Name = "PMTestData"
$pm.Rows("2:5").Delete() | Out-Null

# --- SNMTestData ------------------------------------------------------
# Another full copy of logindata, placed right after PMTestData. This
# one keeps all of its rows (A1:C7) untouched.
$login.Copy($null, $pm) | Out-Null
$snm = $wb.ActiveSheet
$snm.Name = "SNMTestData"

# --- logindata ----------------------------------------------------
# Drop the "create/edit user" rows now that they live in PMTestData.
$login.Rows("6:7").Delete() | Out-Null

# --- Selections / active sheet ---------------------------------------
$pm.Activate()
$pm.Range("D6").Select() | Out-Null

$snm.Activate()
$snm.Range("D7").Select() | Out-Null

$login.Activate()
$login.Range("B18").Select() | Out-Null
